$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hyperlinks first: this both (a) establishes the built-in "Hyperlink"
# cell style as cellXfs index 1 (before the currency numFmt style is
# created further down, which must land on index 2), and (b) creates the
# rIds in the exact order rId1->C10, rId2->C6, rId3->C9. The placeholder
# text Excel auto-fills from the address is overwritten below, so it does
# not end up in the saved shared-string table.
$null = $ws.Hyperlinks.Add($ws.Range("C10"), "https://www.aliexpress.com/wholesale?SearchText=gps+tracker")
$null = $ws.Hyperlinks.Add($ws.Range("C6"), "https://www.aliexpress.com/wholesale?SearchText=temperature+humidity+sensor")
$null = $ws.Hyperlinks.Add($ws.Range("C9"), "https://www.aliexpress.com/wholesale?SearchText=arduino")

# NOTE: cell-write order below is deliberately chosen (not simple row order)
# so the shared-string table is built up in the same sequence the original
# author produced: naam:, link:, sensoren:, licht, CO2, gps tracker,
# kostprijs:, link, 12,78 $, temperatuur & humidity, arduino, 22,00 $.

$ws.Range("A1").Value = "naam:"
$ws.Range("B1").ClearContents()
$ws.Range("C1").Value = "link: "

$ws.Range("A3").Value = "sensoren:"
$ws.Range("A4").Value = "licht"
$ws.Range("A5").Value = "CO²"

$ws.Range("A10").Value = "gps tracker"

$ws.Range("D1").Value = "kostprijs:"

$ws.Range("C10").Value = "link"

$ws.Range("D10").Value = "12,78 $"

$ws.Range("A6").Value = "temperatuur & humidity"

$ws.Range("A9").Value = "arduino"

$ws.Range("D9").Value = "22,00 $"

$ws.Range("C6").Value = "link"
$ws.Range("C9").Value = "link"

# --- Numeric/currency cell (creates cellXfs index 2 = currency numFmt) ----
$ws.Range("D6").Value = 2.13
$ws.Range("D6").NumberFormat = "#,##0.00\ [$€-1];[Red]\-#,##0.00\ [$€-1]"
$ws.Range("D6").HorizontalAlignment = -4131

# --- Column widths -------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 22.5
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 11.666666666666668

# --- Page setup ------------------------------------------------------------
$ws.PageSetup.Orientation = 1

# --- Selection -------------------------------------------------------------
$null = $ws.Range("C9").Select()
